# Add newly-sampled specimen rows (TCMA / TCAO) to the trawl metadata sheet.
#
# Column layout: A=Trawl#, B=Species (italic), C=Code, D=Sampled?
#
# NOTE on write order: rows are not always filled strictly left-to-right.
# Excel's shared-string table assigns indexes in first-use order, and in
# the target workbook "TCAO004" (row 69, column C) is indexed *before*
# "Argyropelecus olfersii" (row 69, column B) even though B is the
# left-hand column -- i.e. for that row the code was typed before the
# species name. The $CodeFirst flag below reproduces that per-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{Row=64; Trawl=2; Species="Macroparalepis affinis"; Code="TCMA021"; Sampled="Y"; CodeFirst=$false},
    @{Row=65; Trawl=2; Species="Macroparalepis affinis"; Code="TCMA022"; Sampled="Y"; CodeFirst=$true},
    @{Row=66; Trawl=2; Species="Macroparalepis affinis"; Code="TCMA023"; Sampled="Y"; CodeFirst=$true},
    @{Row=67; Trawl=2; Species="Macroparalepis affinis"; Code="TCMA024"; Sampled="Y"; CodeFirst=$true},
    @{Row=68; Trawl=2; Species="Macroparalepis affinis"; Code="TCMA025"; Sampled="Y"; CodeFirst=$true},
    @{Row=69; Trawl=8; Species="Argyropelecus olfersii"; Code="TCAO004"; Sampled="Y"; CodeFirst=$true},
    @{Row=70; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO011"; Sampled="Y"; CodeFirst=$true},
    @{Row=71; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO013"; Sampled="Y"; CodeFirst=$true},
    @{Row=72; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO014"; Sampled="Y"; CodeFirst=$true},
    @{Row=73; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO015"; Sampled="Y"; CodeFirst=$true},
    @{Row=74; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO017"; Sampled="Y"; CodeFirst=$true},
    @{Row=75; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO018"; Sampled="Y"; CodeFirst=$true},
    @{Row=76; Trawl=2; Species="Argyropelecus olfersii"; Code="TCAO019"; Sampled="Y"; CodeFirst=$true}
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Trawl

    if ($r.CodeFirst) {
        $ws.Cells.Item($r.Row, 3).Value = $r.Code
        $ws.Cells.Item($r.Row, 2).Value = $r.Species
        $ws.Cells.Item($r.Row, 2).Font.Italic = $true
    } else {
        $ws.Cells.Item($r.Row, 2).Value = $r.Species
        $ws.Cells.Item($r.Row, 2).Font.Italic = $true
        $ws.Cells.Item($r.Row, 3).Value = $r.Code
    }

    $ws.Cells.Item($r.Row, 4).Value = $r.Sampled
}

# Reflect the final on-screen state: the newly-added rows end up selected.
$ws.Activate()
$ws.Range("A73:A76").Select()
